$wb = $excel.ActiveWorkbook

# Sheet "展览" is the 1st worksheet (sheet1.xml)
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F6").Value = 5403
$wsExhibit.Range("F8").Value = 5390
$wsExhibit.Range("F11").Value = 1380

# Sheet "全部类型" is the 4th worksheet (sheet4.xml)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F7").Value = 5403
$wsAll.Range("F9").Value = 5390
$wsAll.Range("F12").Value = 1380
